$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Xpath values to drop the "[1]" index predicate
$ws.Range("B2").Value = "/NewDataSet/Table/Town"
$ws.Range("B3").Value = "/NewDataSet/Table/County"
$ws.Range("B4").Value = "/NewDataSet/Table/PostCode"

# Update the selected range/active cell shown in the sheet view
$ws.Range("F8").Select()
